$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C21").Value = 7040
$ws.Range("C22").Value = 2368
$ws.Range("C23").Value = 2368

$ws.Range("C27").Value = 7042
$ws.Range("C28").Value = 2402
$ws.Range("C29").Value = 2402

$ws.Range("F22").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
